$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (new quarters added)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formats from F:G (old D:E, now shifted) into new D:E columns
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new D/E quarter values
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 285400
$ws.Range("E8").Value = 281700
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 800
$ws.Range("E14").Value = 3300
$ws.Range("D15").Value = 49200
$ws.Range("E15").Value = 45800
$ws.Range("D17").Value = 222000
$ws.Range("E17").Value = 214700
$ws.Range("D18").Value = 63400
$ws.Range("E18").Value = 67000
$ws.Range("D20").Value = -700
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 120700
$ws.Range("E21").Value = 121000
$ws.Range("D22").Value = 12600
$ws.Range("E22").Value = 13200
$ws.Range("D23").Value = 50100
$ws.Range("E23").Value = 53800
$ws.Range("D24").Value = 7300
$ws.Range("E24").Value = 10800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 42800
$ws.Range("E26").Value = 43000
$ws.Range("D27").Value = 42800
$ws.Range("E27").Value = 43000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 700
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 42800
$ws.Range("E33").Value = 43000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 42800
$ws.Range("E35").Value = 43000
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 20300
$ws.Range("E41").Value = 16900
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 193300
$ws.Range("E43").Value = 195100
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 52500
$ws.Range("E45").Value = 48300
$ws.Range("D46").Value = 266100
$ws.Range("E46").Value = 260300
$ws.Range("D47").Value = 22000
$ws.Range("E47").Value = 18700
$ws.Range("D48").Value = 177100
$ws.Range("E48").Value = 176100
$ws.Range("D49").Value = 2923300
$ws.Range("E49").Value = 2907000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 264900
$ws.Range("E52").Value = 281100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 3653400
$ws.Range("E54").Value = 3643200
$ws.Range("D57").Value = 67800
$ws.Range("E57").Value = 43900
$ws.Range("D58").Value = 52500
$ws.Range("E58").Value = 47000
$ws.Range("D59").Value = 118700
$ws.Range("E59").Value = 101000
$ws.Range("D60").Value = 239000
$ws.Range("E60").Value = 191900
$ws.Range("D61").Value = 1284200
$ws.Range("E61").Value = 1351000
$ws.Range("D62").Value = 343700
$ws.Range("E62").Value = 360700
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1866900
$ws.Range("E66").Value = 1903600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 381100
$ws.Range("E72").Value = 338300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 1786500
$ws.Range("E76").Value = 1739600
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 42800
$ws.Range("E81").Value = 43000
$ws.Range("D83").Value = 58000
$ws.Range("E83").Value = 54000
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 130100
$ws.Range("E89").Value = 95400
$ws.Range("D91").Value = -9900
$ws.Range("E91").Value = -5600
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -64400
$ws.Range("E94").Value = -21700
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -62300
$ws.Range("E100").Value = -68200
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 3400
$ws.Range("E102").Value = 5500

# Row 91 (Capital Expenditures) received corrected historical data beyond the shift
$ws.Range("F91").Value = -4400
$ws.Range("G91").Value = -10100
$ws.Range("H91").Value = -22100
$ws.Range("I91").Value = -1300
$ws.Range("J91").Value = -300
